$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "county_name"
$ws.Range("B1").Value = "fips"
$ws.Range("C1").Value = "level1_dum"
$ws.Range("D1").Value = "total_level1"
$ws.Range("E1").Value = "level2_dum"
$ws.Range("F1").Value = "total_level2"
$ws.Range("G1").Value = "pop_2010"
$ws.Range("H1").Value = "rucc_2013"
$ws.Range("I1").Value = "metro_non"

# --- Data rows 2-65 ---
$countyNames = @(
  @("Adams", 8001, 0, 0, 0, 0, 441603, 1, 1),
  @("Alamosa", 8003, 0, 0, 0, 0, 15445, 7, 0),
  @("Arapahoe", 8005, 1, 1, 1, 2, 572003, 1, 1),
  @("Archuleta", 8007, 0, 0, 0, 0, 12084, 7, 0),
  @("Baca", 8009, 0, 0, 0, 0, 3788, 9, 0),
  @("Bent", 8011, 0, 0, 0, 0, 6499, 7, 0),
  @("Broomfield", 8014, 0, 0, 0, 0, 55889, 1, 1),
  @("Boulder", 8013, 0, 0, 1, 2, 294567, 2, 1),
  @("Chaffee", 8015, 0, 0, 0, 0, 17809, 7, 0),
  @("Cheyenne", 8017, 0, 0, 0, 0, 1836, 9, 0),
  @("Clear Creek", 8019, 0, 0, 0, 0, 9088, 1, 1),
  @("Conejos", 8021, 0, 0, 0, 0, 8256, 9, 0),
  @("Costilla", 8023, 0, 0, 0, 0, 3524, 9, 0),
  @("Crowley", 8025, 0, 0, 0, 0, 5823, 8, 0),
  @("Custer", 8027, 0, 0, 0, 0, 4255, 8, 0),
  @("Delta", 8029, 0, 0, 0, 0, 30952, 6, 0),
  @("Denver", 8031, 1, 1, 0, 0, 600158, 1, 1),
  @("Dolores", 8033, 0, 0, 0, 0, 2064, 9, 0),
  @("Douglas", 8035, 0, 0, 1, 2, 285465, 1, 1),
  @("Eagle", 8037, 0, 0, 0, 0, 52197, 5, 0),
  @("Elbert", 8039, 0, 0, 0, 0, 23086, 1, 1),
  @("El Paso", 8041, 1, 1, 1, 1, 622263, 2, 1),
  @("Fremont", 8043, 0, 0, 0, 0, 46824, 4, 0),
  @("Garfield", 8045, 0, 0, 0, 0, 56389, 5, 0),
  @("Gilpin", 8047, 0, 0, 0, 0, 5441, 1, 1),
  @("Grand", 8049, 0, 0, 0, 0, 14843, 7, 0),
  @("Gunnison", 8051, 0, 0, 0, 0, 15324, 7, 0),
  @("Hinsdale", 8053, 0, 0, 0, 0, 843, 9, 0),
  @("Huerfano", 8055, 0, 0, 0, 0, 6711, 6, 0),
  @("Jackson", 8057, 0, 0, 0, 0, 1394, 9, 0),
  @("Jefferson", 8059, 1, 1, 0, 0, 534543, 1, 1),
  @("Kiowa", 8061, 0, 0, 0, 0, 1398, 9, 0),
  @("Kit Carson", 8063, 0, 0, 0, 0, 8270, 7, 0),
  @("Lake", 8065, 0, 0, 0, 0, 7310, 6, 0),
  @("La Plata", 8067, 0, 0, 0, 0, 51334, 4, 0),
  @("Larimer", 8069, 0, 0, 1, 1, 299630, 2, 1),
  @("Las Animas", 8071, 0, 0, 0, 0, 15507, 7, 0),
  @("Lincoln", 8073, 0, 0, 0, 0, 5467, 8, 0),
  @("Logan", 8075, 0, 0, 0, 0, 22709, 7, 0),
  @("Mesa", 8077, 0, 0, 1, 1, 146723, 3, 1),
  @("Mineral", 8079, 0, 0, 0, 0, 712, 9, 0),
  @("Moffat", 8081, 0, 0, 0, 0, 13795, 7, 0),
  @("Montezuma", 8083, 0, 0, 0, 0, 25535, 6, 0),
  @("Montrose", 8085, 0, 0, 0, 0, 41276, 4, 0),
  @("Morgan", 8087, 0, 0, 0, 0, 28159, 6, 0),
  @("Otero", 8089, 0, 0, 0, 0, 18831, 6, 0),
  @("Ouray", 8091, 0, 0, 0, 0, 4436, 9, 0),
  @("Park", 8093, 0, 0, 0, 0, 16206, 1, 1),
  @("Phillips", 8095, 0, 0, 0, 0, 4442, 9, 0),
  @("Pitkin", 8097, 0, 0, 0, 0, 17148, 7, 0),
  @("Prowers", 8099, 0, 0, 0, 0, 12551, 7, 0),
  @("Pueblo", 8101, 0, 0, 1, 1, 159063, 3, 1),
  @("Rio Blanco", 8103, 0, 0, 0, 0, 6666, 9, 0),
  @("Rio Grande", 8105, 0, 0, 0, 0, 11982, 7, 0),
  @("Routt", 8107, 0, 0, 0, 0, 23509, 7, 0),
  @("Saguache", 8109, 0, 0, 0, 0, 6108, 9, 0),
  @("San Juan", 8111, 0, 0, 0, 0, 699, 9, 0),
  @("San Miguel", 8113, 0, 0, 0, 0, 7359, 9, 0),
  @("Sedgwick", 8115, 0, 0, 0, 0, 2379, 9, 0),
  @("Summit", 8117, 0, 0, 0, 0, 27994, 5, 0),
  @("Teller", 8119, 0, 0, 0, 0, 23350, 2, 1),
  @("Washington", 8121, 0, 0, 0, 0, 4814, 9, 0),
  @("Weld", 8123, 0, 0, 1, 1, 252825, 2, 1),
  @("Yuma", 8125, 0, 0, 0, 0, 10043, 7, 0)
)

$r = 2
foreach ($row in $countyNames) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

# --- Stray formatted (but empty) cells left over from the author's session ---
$ws.Range("A2").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- View state: zoom + active selection ---
$null = ($ws.Application.ActiveWindow.Zoom = 125)
$ws.Range("J8").Select() | Out-Null